$d = $word.ActiveDocument

$pairs = @(
    @("388÷7=", "881÷2="),
    @("265÷8=", "148÷5="),
    @("383÷9=", "436÷5="),
    @("162÷8=", "640÷4="),
    @("733÷5=", "918÷4="),
    @("153÷6=", "652÷4="),
    @("322÷9=", "915÷2="),
    @("233÷4=", "494÷3="),
    @("809÷5=", "320÷3="),
    @("400÷9=", "212÷7="),
    @("558÷5=", "901÷5="),
    @("936÷9=", "672÷8="),
    @("198÷2=", "361÷9="),
    @("228÷7=", "538÷2="),
    @("483÷7=", "700÷6="),
    @("525÷3=", "302÷4="),
    @("306÷9=", "289÷5="),
    @("216÷9=", "487÷7="),
    @("865÷9=", "201÷8="),
    @("464÷5=", "978÷9="),
    @("427÷6=", "379÷9="),
    @("254÷5=", "388÷5="),
    @("977÷6=", "733÷9="),
    @("485÷3=", "771÷7="),
    @("479÷4=", "717÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
